$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.2341203333333333
$ws.Range("H2").Value = 0.702361
$ws.Range("I2").Value = 0.1153892343949327
$ws.Range("J2").Value = 0.1153892343949327
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 92.59233966666666
$ws.Range("N2").Value = 277.777019
$ws.Range("O2").Value = 0.5465415277631132
$ws.Range("P2").Value = 0.5465415277631133
$ws.Range("Q2").Value = 21.67774942687322
$ws.Range("R2").Value = 195.099744841859
$ws.Range("S2").Value = 0.06306500845362249
$ws.Range("T2").Value = 0.0630650084536225

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.2341203333333333
$ws.Range("H3").Value = 0.702361
$ws.Range("I3").Value = 0.1153892343949327
$ws.Range("J3").Value = 0.1153892343949327
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 15.14173
$ws.Range("N3").Value = 45.42519
$ws.Range("O3").Value = 0.08937655401050183
$ws.Range("P3").Value = 0.08937655401050183
$ws.Range("Q3").Value = 3.544986874843334
$ws.Range("R3").Value = 31.90488187359
$ws.Range("S3").Value = 0.01031309214012916
$ws.Range("T3").Value = 0.01031309214012916

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.2341203333333333
$ws.Range("H4").Value = 0.702361
$ws.Range("I4").Value = 0.1153892343949327
$ws.Range("J4").Value = 0.1153892343949327
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 28.17812733333333
$ws.Range("N4").Value = 84.534382
$ws.Range("O4").Value = 0.1663260353686444
$ws.Range("P4").Value = 0.1663260353686444
$ws.Range("Q4").Value = 6.597072563989111
$ws.Range("R4").Value = 59.373653075902
$ws.Range("S4").Value = 0.01919223388113238
$ws.Range("T4").Value = 0.01919223388113238

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.2341203333333333
$ws.Range("H5").Value = 0.702361
$ws.Range("I5").Value = 0.1153892343949327
$ws.Range("J5").Value = 0.1153892343949327
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 33.50281533333334
$ws.Range("N5").Value = 100.508446
$ws.Range("O5").Value = 0.1977558828577406
$ws.Range("P5").Value = 0.1977558828577405
$ws.Range("Q5").Value = 7.843690293445112
$ws.Range("R5").Value = 70.593212641006
$ws.Range("S5").Value = 0.02281889992004868
$ws.Range("T5").Value = 0.02281889992004868

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.794841333333333
$ws.Range("H6").Value = 5.384524
$ws.Range("I6").Value = 0.8846107656050674
$ws.Range("J6").Value = 0.8846107656050674
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 92.59233966666666
$ws.Range("N6").Value = 277.777019
$ws.Range("O6").Value = 0.5465415277631132
$ws.Range("P6").Value = 0.5465415277631133
$ws.Range("Q6").Value = 166.1885583837729
$ws.Range("R6").Value = 1495.697025453956
$ws.Range("S6").Value = 0.4834765193094907
$ws.Range("T6").Value = 0.4834765193094908

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.794841333333333
$ws.Range("H7").Value = 5.384524
$ws.Range("I7").Value = 0.8846107656050674
$ws.Range("J7").Value = 0.8846107656050674
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 15.14173
$ws.Range("N7").Value = 45.42519
$ws.Range("O7").Value = 0.08937655401050183
$ws.Range("P7").Value = 0.08937655401050183
$ws.Range("Q7").Value = 27.17700286217334
$ws.Range("R7").Value = 244.59302575956
$ws.Range("S7").Value = 0.07906346187037268
$ws.Range("T7").Value = 0.07906346187037268

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.794841333333333
$ws.Range("H8").Value = 5.384524
$ws.Range("I8").Value = 0.8846107656050674
$ws.Range("J8").Value = 0.8846107656050674
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 28.17812733333333
$ws.Range("N8").Value = 84.534382
$ws.Range("O8").Value = 0.1663260353686444
$ws.Range("P8").Value = 0.1663260353686444
$ws.Range("Q8").Value = 50.57526763379644
$ws.Range("R8").Value = 455.1774087041679
$ws.Range("S8").Value = 0.1471338014875121
$ws.Range("T8").Value = 0.1471338014875121

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.794841333333333
$ws.Range("H9").Value = 5.384524
$ws.Range("I9").Value = 0.8846107656050674
$ws.Range("J9").Value = 0.8846107656050674
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 33.50281533333334
$ws.Range("N9").Value = 100.508446
$ws.Range("O9").Value = 0.1977558828577406
$ws.Range("P9").Value = 0.1977558828577405
$ws.Range("Q9").Value = 60.13223774330045
$ws.Range("R9").Value = 541.190139689704
$ws.Range("S9").Value = 0.1749369829376919
$ws.Range("T9").Value = 0.1749369829376919

